$wb = $excel.ActiveWorkbook

# Update worksheet cell values for the 2025-11-24 data refresh.
# Each block below targets one worksheet, setting the numeric cell values
# to match the updated source data (year-2025 "L" column reflects the new
# year-to-date totals as of 2025-11-24; a couple of 2016 "C" column values
# were also corrected).

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 6000
$ws.Range("L3").Value = 6545
$ws.Range("C4").Value = 1876
$ws.Range("L4").Value = 1608
$ws.Range("L5").Value = 390
$ws.Range("L6").Value = 5368
$ws.Range("C7").Value = 28420
$ws.Range("L7").Value = 19911

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L2").Value = 66
$ws.Range("L7").Value = 220

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 396
$ws.Range("L3").Value = 465
$ws.Range("L4").Value = 91
$ws.Range("L7").Value = 1319

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L3").Value = 175
$ws.Range("L6").Value = 95
$ws.Range("L7").Value = 433

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L3").Value = 314
$ws.Range("L7").Value = 895

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L2").Value = 103
$ws.Range("L3").Value = 91
$ws.Range("L7").Value = 285

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 228
$ws.Range("L3").Value = 270
$ws.Range("L4").Value = 40
$ws.Range("L7").Value = 759

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 143
$ws.Range("L3").Value = 128
$ws.Range("L7").Value = 391

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L2").Value = 100
$ws.Range("L3").Value = 141
$ws.Range("L7").Value = 346

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L4").Value = 70
$ws.Range("L7").Value = 643
$ws.Range("L8").Value = 1319
$ws.Range("L9").Value = 114
$ws.Range("L14").Value = 98
$ws.Range("L16").Value = 42
$ws.Range("L19").Value = 538
$ws.Range("L20").Value = 506
$ws.Range("L21").Value = 62
$ws.Range("L23").Value = 213
$ws.Range("L29").Value = 1119
$ws.Range("L30").Value = 88
$ws.Range("L33").Value = 895
$ws.Range("L37").Value = 759
$ws.Range("L42").Value = 639
$ws.Range("L47").Value = 140
$ws.Range("L52").Value = 422
$ws.Range("L53").Value = 220
$ws.Range("L54").Value = 434
$ws.Range("L55").Value = 209
$ws.Range("C63").Value = 300
$ws.Range("L63").Value = 57
$ws.Range("L64").Value = 127
$ws.Range("L65").Value = 391
$ws.Range("L67").Value = 687
$ws.Range("L76").Value = 303
$ws.Range("L79").Value = 554
$ws.Range("L83").Value = 433
$ws.Range("L85").Value = 985
$ws.Range("L89").Value = 277
$ws.Range("L91").Value = 268
$ws.Range("L94").Value = 249
$ws.Range("L95").Value = 285
$ws.Range("L96").Value = 222
$ws.Range("L98").Value = 107
$ws.Range("L99").Value = 346
$ws.Range("L100").Value = 38
$ws.Range("C101").Value = 28420
$ws.Range("L101").Value = 19911

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 193
$ws.Range("L3").Value = 269
$ws.Range("L5").Value = 20
$ws.Range("L7").Value = 687

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L5").Value = 4
$ws.Range("L7").Value = 434

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 331
$ws.Range("L6").Value = 274
$ws.Range("L7").Value = 1119

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 191
$ws.Range("L7").Value = 538

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L6").Value = 135
$ws.Range("L7").Value = 303

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("L2").Value = 40
$ws.Range("L7").Value = 98

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 174
$ws.Range("L3").Value = 221
$ws.Range("L6").Value = 177
$ws.Range("L7").Value = 639

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L3").Value = 70
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 209

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L2").Value = 55
$ws.Range("L7").Value = 213

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L2").Value = 70
$ws.Range("L7").Value = 222

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L2").Value = 92
$ws.Range("L4").Value = 14
$ws.Range("L7").Value = 268

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("L3").Value = 15
$ws.Range("L7").Value = 62

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 173
$ws.Range("L6").Value = 150
$ws.Range("L7").Value = 554

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L3").Value = 37
$ws.Range("L6").Value = 35
$ws.Range("L7").Value = 127

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 159
$ws.Range("L3").Value = 177
$ws.Range("L7").Value = 506

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 38

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L6").Value = 155
$ws.Range("L7").Value = 643

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L4").Value = 31
$ws.Range("L7").Value = 249

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L3").Value = 47
$ws.Range("L4").Value = 10
$ws.Range("L7").Value = 140

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("L3").Value = 18
$ws.Range("L4").Value = 12
$ws.Range("L7").Value = 107

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("L2").Value = 34
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L2").Value = 70
$ws.Range("L7").Value = 277

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 295
$ws.Range("L3").Value = 409
$ws.Range("L7").Value = 985

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L3").Value = 133
$ws.Range("L7").Value = 422

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 42
